# Remove the "NA"/"na" placeholder entries, leaving the cells blank instead.
# (commit message: "removed NA's and left blanks instead.")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose content is the literal "NA" (rows 2-9) or "na" (rows 10-17).
# ClearContents() (rather than setting .Value = "") avoids introducing
# incidental row-height attributes on the previously-hidden rows.
$naCells = @("C2", "C3", "B4", "C4", "B9", "C9", "C10", "C11", "B12", "C12", "B17", "C17")
foreach ($addr in $naCells) {
    $ws.Range($addr).ClearContents()
}

# The sheet had an AutoFilter on column H restricted to "A" (Workout A),
# which hid rows 10-17. Showing all data drops that filter criterion and
# unhides the rows, while keeping the AutoFilter dropdowns in place.
$ws.ShowAllData()

# Match the saved selection left behind by the edit.
$ws.Range("E21").Select()
